$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.8710206747055054
$ws.Range("B1").Value = 2.262847185134888
$ws.Range("C1").Value = 1.577043414115906
$ws.Range("D1").Value = 1.386647820472717
$ws.Range("E1").Value = 1.464355111122131
